$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (A10)
$ws.Range("A10").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-07-08 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-7
$ws.Range("D2").Value = 0.2636988995744722
$ws.Range("E2").Value = -0.008115143163374694

$ws.Range("D3").Value = 0.5300434519356292
$ws.Range("E3").Value = -0.01135763338615947

$ws.Range("D4").Value = 0.05259768506354875
$ws.Range("E4").Value = -0.008733624454148603

$ws.Range("D5").Value = 0.09617047197991686
$ws.Range("E5").Value = -0.01546534106600383

$ws.Range("D6").Value = 0.057489491446433
$ws.Range("E6").Value = -0.01035103510351043

$ws.Range("E7").Value = -0.01070174685008607

# Restore sheet protection that was temporarily lifted to allow the edits above
$ws.Protect()
